$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (want-to-go count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 913
$wsExhibit.Range("F5").Value = 538

# Sheet "全部类型" (All Types) - same two events appear here, update accordingly
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 913
$wsAll.Range("F6").Value = 538
